$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Find the last used row in column A (the "Beteckning" column) to know how far data extends.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Update the "Förändrad" column (C) date serial value from 45184 to 45186 for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# Update the HYPERLINK formulas on row 2 to include a friendly display text
# (the "Beteckning" value of that row) as the second HYPERLINK argument.
$label = $ws.Cells.Item(2, 1).Value2

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/artfynd/A 16408-2019.xlsx", "' + $label + '")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/kartor/A 16408-2019.png", "' + $label + '")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/klagomål/A 16408-2019.docx", "' + $label + '")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/klagomålsmail/A 16408-2019.docx", "' + $label + '")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/tillsyn/A 16408-2019.docx", "' + $label + '")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_ESSUNGA/tillsynsmail/A 16408-2019.docx", "' + $label + '")'
